$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Aarhus University
$ws.Range("C4").Value = 157
$ws.Range("D4").Value = 84.40000000000001
$ws.Range("E4").Value = 78.5
$ws.Range("F4").Value = 88.90000000000001

# Row 5 - Aarhus University Hospital
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 78.09999999999999
$ws.Range("E5").Value = 66.60000000000001
$ws.Range("F5").Value = 86.5

# Row 8 - Copenhagen University Hospital
$ws.Range("C8").Value = 138
$ws.Range("D8").Value = 80.2
$ws.Range("E8").Value = 73.59999999999999
$ws.Range("F8").Value = 85.5

# Row 16 - Karolinska Institutet
$ws.Range("C16").Value = 132
$ws.Range("D16").Value = 79.5
$ws.Range("E16").Value = 72.7
$ws.Range("F16").Value = 85

# Row 26 - Odense University Hospital
$ws.Range("C26").Value = 81
$ws.Range("D26").Value = 89
$ws.Range("E26").Value = 80.90000000000001
$ws.Range("F26").Value = 93.89999999999999

# Row 44 - University of Copenhagen
$ws.Range("C44").Value = 72
$ws.Range("D44").Value = 74.2
$ws.Range("E44").Value = 64.7
$ws.Range("F44").Value = 81.89999999999999

# Row 49 - University of Oulu
$ws.Range("C49").Value = 18
$ws.Range("D49").Value = 72
$ws.Range("E49").Value = 52.40000000000001
$ws.Range("F49").Value = 85.7
